$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 96, pushing existing rows 96-101 down to 97-102.
$ws.Rows.Item(96).Insert()

# Populate the new weekly record in row 96.
$ws.Range("A96").Value = 1
$ws.Range("B96").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C96").Value = "Arica y Parinacota"
$ws.Range("D96").Value = 44714
$ws.Range("E96").Value = 15
$ws.Range("F96").Value = "Fruta"
$ws.Range("G96").Value = 100102
$ws.Range("H96").Value = "Cítricos"
$ws.Range("I96").Value = 100102004
$ws.Range("J96").Value = "Mandarina"
$ws.Range("K96").Value = "Clemenuless"
$ws.Range("L96").Value = "Segunda"
$ws.Range("M96").Value = 300
$ws.Range("N96").Value = 18000
$ws.Range("O96").Value = 19000
$ws.Range("P96").Value = 18500
$ws.Range("Q96").Value = "$/caja 20 kilos"
$ws.Range("R96").Value = "Región de Coquimbo"
$ws.Range("S96").Value = 925
$ws.Range("T96").Value = 20
